$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.155.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.313.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.93%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.310.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.32%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "

# Row 12
$ws.Range("E12").Value = "  -0.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.966.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.36%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.84%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.720.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.311.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.92%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.49%  "

# Row 23
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.53%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.80%  "

# Row 29
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

# Row 31
$ws.Range("E31").Value = "  -3.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0719"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.68%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.70%  "

# Row 34
$ws.Range("E34").Value = "  -6.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.380"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.29%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "

# Row 38
$ws.Range("E38").Value = "  +0.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.87%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "317.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.567"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0490"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0213"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50
$ws.Range("E50").Value = "  +1.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "
